$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4854.8066
$ws.Range("J64").Value = 4931
$ws.Range("L64").Value = 4931
$ws.Range("N64").Value = -5427
$ws.Range("H67").Value = 4854.8066
$ws.Range("J67").Value = 4931
$ws.Range("L67").Value = 4931
$ws.Range("N67").Value = -6647
$ws.Range("H103").Value = 535.1429000000001
$ws.Range("J103").Value = 550
$ws.Range("L103").Value = 1650
$ws.Range("N103").Value = -2822
$ws.Range("H113").Value = 35717556
$ws.Range("J113").Value = 50003876
$ws.Range("L113").Value = 50003876
$ws.Range("N113").Value = -50010384
$ws.Range("H138").Value = 2330.2654
$ws.Range("I138").Value = 1146
$ws.Range("J138").Value = 2561.3416
$ws.Range("K138").Value = 3438
$ws.Range("L138").Value = 7684.024800000001
$ws.Range("M138").Value = 1702
$ws.Range("N138").Value = -17964.0248
$ws.Range("H139").Value = 72000
$ws.Range("J139").Value = 72000
$ws.Range("L139").Value = 72000
$ws.Range("N139").Value = -82280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10036173
$ws.Range("I32").Value = 14746299
$ws.Range("J32").Value = 27156.625
$ws.Range("K32").Value = 14746299
$ws.Range("L32").Value = 27156.625
$ws.Range("M32").Value = -14746012
$ws.Range("N32").Value = -27730.625
$ws.Range("H36").Value = 6499.5
$ws.Range("I36").Value = 6499.5
$ws.Range("K36").Value = 6499.5
$ws.Range("M36").Value = -6153.5
$ws.Range("H45").Value = 100004184
$ws.Range("I45").Value = 166671330
$ws.Range("J45").Value = 3457
$ws.Range("K45").Value = 166671330
$ws.Range("L45").Value = 3457
$ws.Range("M45").Value = -166670953
$ws.Range("N45").Value = -4211
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H122").Value = 3122.3845
$ws.Range("I122").Value = 2605.1035
$ws.Range("J122").Value = 4622.5
$ws.Range("K122").Value = 7815.310500000001
$ws.Range("L122").Value = 13867.5
$ws.Range("M122").Value = -5365.310500000001
$ws.Range("N122").Value = -18767.5
$ws.Range("H123").Value = 87996.664
$ws.Range("J123").Value = 87996.664
$ws.Range("L123").Value = 87996.664
$ws.Range("N123").Value = -97796.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2313.3333
$ws.Range("J64").Value = 1776
$ws.Range("L64").Value = 1776
$ws.Range("N64").Value = -2226
$ws.Range("H67").Value = 2313.3333
$ws.Range("J67").Value = 1776
$ws.Range("L67").Value = 1776
$ws.Range("N67").Value = -3336
$ws.Range("H99").Value = 7687.6
$ws.Range("I99").Value = 11353
$ws.Range("J99").Value = 4022.2
$ws.Range("K99").Value = 11353
$ws.Range("L99").Value = 4022.2
$ws.Range("M99").Value = -9855
$ws.Range("N99").Value = -7018.2
$ws.Range("H123").Value = 59000
$ws.Range("J123").Value = 59000
$ws.Range("L123").Value = 59000
$ws.Range("N123").Value = -68800
$ws.Range("H125").Value = 62000
$ws.Range("J125").Value = 62000
$ws.Range("L125").Value = 62000
$ws.Range("N125").Value = -71840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6427.25
$ws.Range("I58").Value = 7603.3335
$ws.Range("J58").Value = 2899
$ws.Range("K58").Value = 7603.3335
$ws.Range("L58").Value = 2899
$ws.Range("M58").Value = -7400.3335
$ws.Range("N58").Value = -3305
$ws.Range("H99").Value = 1978.4
$ws.Range("I99").Value = 1858.2858
$ws.Range("K99").Value = 1858.2858
$ws.Range("M99").Value = -360.2858000000001
$ws.Range("H122").Value = 1974.8667
$ws.Range("I122").Value = 1804.6
$ws.Range("J122").Value = 2315.4
$ws.Range("K122").Value = 5413.799999999999
$ws.Range("L122").Value = 6946.200000000001
$ws.Range("M122").Value = -2963.799999999999
$ws.Range("N122").Value = -11846.2
$ws.Range("H126").Value = 1978.4
$ws.Range("I126").Value = 1858.2858
$ws.Range("K126").Value = 5574.857400000001
$ws.Range("M126").Value = -3104.857400000001
$ws.Range("H136").Value = 6427.25
$ws.Range("I136").Value = 7603.3335
$ws.Range("J136").Value = 2899
$ws.Range("K136").Value = 22810.0005
$ws.Range("L136").Value = 8697
$ws.Range("M136").Value = -20260.0005
$ws.Range("N136").Value = -13797

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2021.2273
$ws.Range("J5").Value = 2125.75
$ws.Range("L5").Value = 6377.25
$ws.Range("N5").Value = -6601.25
$ws.Range("H25").Value = 160.66667
$ws.Range("I25").Value = 130
$ws.Range("J25").Value = 222
$ws.Range("K25").Value = 390
$ws.Range("L25").Value = 666
$ws.Range("M25").Value = -221
$ws.Range("N25").Value = -1004
$ws.Range("H30").Value = 160.66667
$ws.Range("I30").Value = 130
$ws.Range("J30").Value = 222
$ws.Range("K30").Value = 390
$ws.Range("L30").Value = 666
$ws.Range("M30").Value = -288
$ws.Range("N30").Value = -870
$ws.Range("H39").Value = 1890.1818
$ws.Range("J39").Value = 2018.6
$ws.Range("L39").Value = 6055.799999999999
$ws.Range("N39").Value = -6643.799999999999
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H132").Value = 2714.1333
$ws.Range("I132").Value = 2333
$ws.Range("J132").Value = 2809.4167
$ws.Range("K132").Value = 20997
$ws.Range("L132").Value = 25284.7503
$ws.Range("M132").Value = -18467
$ws.Range("N132").Value = -30344.7503
$ws.Range("H133").Value = 4166.6665
$ws.Range("I133").Value = 4166.6665
$ws.Range("K133").Value = 12499.9995
$ws.Range("M133").Value = -7439.999500000002
$ws.Range("H135").Value = 2021.2273
$ws.Range("J135").Value = 2125.75
$ws.Range("L135").Value = 19131.75
$ws.Range("N135").Value = -24201.75
$ws.Range("H137").Value = 4243.706
$ws.Range("J137").Value = 5873.75
$ws.Range("L137").Value = 17621.25
$ws.Range("N137").Value = -27821.25
$ws.Range("H139").Value = 3779.8572
$ws.Range("I139").Value = 3779.8572
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 11339.5716
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -6199.571599999999
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 84241.375
$ws.Range("I140").Value = 161491.11
$ws.Range("J140").Value = 2700
$ws.Range("K140").Value = 484473.33
$ws.Range("L140").Value = 8100
$ws.Range("M140").Value = -479293.33
$ws.Range("N140").Value = -18460
$ws.Range("H141").Value = 280720.72
$ws.Range("I141").Value = 339769.78
$ws.Range("K141").Value = 1019309.34
$ws.Range("M141").Value = -1014129.34

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2274.111
$ws.Range("I122").Value = 2274.111
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6822.333
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4372.333
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 4246.231
$ws.Range("I126").Value = 3957.2856
$ws.Range("K126").Value = 11871.8568
$ws.Range("M126").Value = -9401.856800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1226.5
$ws.Range("H27").Value = 1226.5
$ws.Range("H100").Value = 2799.2632
$ws.Range("I100").Value = 2283.5386
$ws.Range("J100").Value = 3916.6667
$ws.Range("K100").Value = 2283.5386
$ws.Range("L100").Value = 3916.6667
$ws.Range("M100").Value = -1742.5386
$ws.Range("N100").Value = -4998.6667
$ws.Range("H115").Value = 133328
$ws.Range("J115").Value = 133328
$ws.Range("L115").Value = 133328
$ws.Range("N115").Value = -135678
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H122").Value = 10514.652
$ws.Range("I122").Value = 9167.117
$ws.Range("J122").Value = 14332.667
$ws.Range("K122").Value = 27501.351
$ws.Range("L122").Value = 42998.001
$ws.Range("M122").Value = -25051.351
$ws.Range("N122").Value = -47898.001
